# Add 5 new "template" comparison columns (AC:AG) mirroring the existing parameter block
# (see column Z) but swapping in different templateName values: mask_01, mask_02, mask_03,
# mask_04, mask_05 (mask_03 repeats the existing run). Filled one column at a time, top to
# bottom, so the column order below matches how the data was entered.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column AC (29)
$ACValues = @(
    @{ Row = 2; Text = "`n" },
    @{ Row = 3; Text = " varifold`n" },
    @{ Row = 4; Text = " 10`n" },
    @{ Row = 5; Text = " 10`n" },
    @{ Row = 6; Text = " 0.05`n" },
    @{ Row = 7; Text = " 11`n" },
    @{ Row = 8; Text = " mask_01.vtk`n" },
    @{ Row = 9; Text = "`n" },
    @{ Row = 10; Text = " `n" },
    @{ Row = 11; Text = " 200.0`n" },
    @{ Row = 12; Text = " 1.3926833677141874`n" },
    @{ Row = 13; Text = " 1.0314726918915678`n" },
    @{ Row = 14; Text = " 0.26489787688798044`n" },
    @{ Row = 15; Text = " 0.7600476611763216`n" },
    @{ Row = 16; Text = " 1.0277075828332634`n" },
    @{ Row = 17; Text = " 1.772498106962496`n" },
    @{ Row = 18; Text = " 5.33055470049103`n" },
    @{ Row = 19; Text = "`n" },
    @{ Row = 20; Text = " `n" },
    @{ Row = 21; Text = " 20.272240183799003`n" },
    @{ Row = 22; Text = " 71.62598275812971`n" },
    @{ Row = 23; Text = " 94.10247184069144`n" },
    @{ Row = 24; Text = " 69.1025917819341`n" },
    @{ Row = 25; Text = " 68.15191786583149`n" },
    @{ Row = 26; Text = " 64.65104088607714`n" },
    @{ Row = 27; Text = "`n" }
)
foreach ($item in $ACValues) {
    $cell = $ws.Cells.Item($item.Row, 29)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Text
}

# Column AD (30)
$ADValues = @(
    @{ Row = 2; Text = "`n" },
    @{ Row = 3; Text = " varifold`n" },
    @{ Row = 4; Text = " 10`n" },
    @{ Row = 5; Text = " 10`n" },
    @{ Row = 6; Text = " 0.05`n" },
    @{ Row = 7; Text = " 11`n" },
    @{ Row = 8; Text = " mask_02.vtk`n" },
    @{ Row = 9; Text = "`n" },
    @{ Row = 10; Text = " `n" },
    @{ Row = 11; Text = " 200.0`n" },
    @{ Row = 12; Text = " 1.4373899324715347`n" },
    @{ Row = 13; Text = " 1.1170845689556579`n" },
    @{ Row = 14; Text = " 0.32112101507938834`n" },
    @{ Row = 15; Text = " 0.6902016619337479`n" },
    @{ Row = 16; Text = " 1.1030417905326813`n" },
    @{ Row = 17; Text = " 1.6429802035152223`n" },
    @{ Row = 18; Text = " 6.431407656020065`n" },
    @{ Row = 19; Text = "`n" },
    @{ Row = 20; Text = " `n" },
    @{ Row = 21; Text = " 72.64165987741085`n" },
    @{ Row = 22; Text = " 19.79589315624266`n" },
    @{ Row = 23; Text = " 93.30029699701637`n" },
    @{ Row = 24; Text = " 56.11929817798953`n" },
    @{ Row = 25; Text = " 62.38243897133723`n" },
    @{ Row = 26; Text = " 60.84791743599933`n" },
    @{ Row = 27; Text = "`n" }
)
foreach ($item in $ADValues) {
    $cell = $ws.Cells.Item($item.Row, 30)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Text
}

# Column AE (31)
$AEValues = @(
    @{ Row = 2; Text = "`n" },
    @{ Row = 3; Text = " varifold`n" },
    @{ Row = 4; Text = " 10`n" },
    @{ Row = 5; Text = " 10`n" },
    @{ Row = 6; Text = " 0.05`n" },
    @{ Row = 7; Text = " 11`n" },
    @{ Row = 8; Text = " mask_03.vtk`n" },
    @{ Row = 9; Text = "`n" },
    @{ Row = 10; Text = " `n" },
    @{ Row = 11; Text = " 200.0`n" },
    @{ Row = 12; Text = " 1.7568189754609869`n" },
    @{ Row = 13; Text = " 1.3666386634951826`n" },
    @{ Row = 14; Text = " 0.17889475616575568`n" },
    @{ Row = 15; Text = " 0.8497399383398375`n" },
    @{ Row = 16; Text = " 1.3178937037555252`n" },
    @{ Row = 17; Text = " 2.0162563202574564`n" },
    @{ Row = 18; Text = " 7.777727032876133`n" },
    @{ Row = 19; Text = "`n" },
    @{ Row = 20; Text = " `n" },
    @{ Row = 21; Text = " 58.99528015720691`n" },
    @{ Row = 22; Text = " 65.38610769537314`n" },
    @{ Row = 23; Text = " 19.907266531468363`n" },
    @{ Row = 24; Text = " 58.1791096975447`n" },
    @{ Row = 25; Text = " 55.409472574152474`n" },
    @{ Row = 26; Text = " 51.57544733114912`n" },
    @{ Row = 27; Text = "`n" }
)
foreach ($item in $AEValues) {
    $cell = $ws.Cells.Item($item.Row, 31)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Text
}

# Column AF (32)
$AFValues = @(
    @{ Row = 2; Text = "`n" },
    @{ Row = 3; Text = " varifold`n" },
    @{ Row = 4; Text = " 10`n" },
    @{ Row = 5; Text = " 10`n" },
    @{ Row = 6; Text = " 0.05`n" },
    @{ Row = 7; Text = " 11`n" },
    @{ Row = 8; Text = " mask_04.vtk`n" },
    @{ Row = 9; Text = "`n" },
    @{ Row = 10; Text = " `n" },
    @{ Row = 11; Text = " 200.0`n" },
    @{ Row = 12; Text = " 1.43388334231223`n" },
    @{ Row = 13; Text = " 1.0149048194790486`n" },
    @{ Row = 14; Text = " 0.3236543978383953`n" },
    @{ Row = 15; Text = " 0.7662460558113806`n" },
    @{ Row = 16; Text = " 1.144842288192387`n" },
    @{ Row = 17; Text = " 1.6667852167232455`n" },
    @{ Row = 18; Text = " 6.0562791602653085`n" },
    @{ Row = 19; Text = "`n" },
    @{ Row = 20; Text = " `n" },
    @{ Row = 21; Text = " 65.2162525193901`n" },
    @{ Row = 22; Text = " 60.094306984538576`n" },
    @{ Row = 23; Text = " 95.30408078022428`n" },
    @{ Row = 24; Text = " 16.22260601495686`n" },
    @{ Row = 25; Text = " 65.94857783792861`n" },
    @{ Row = 26; Text = " 60.557164827407675`n" },
    @{ Row = 27; Text = "`n" }
)
foreach ($item in $AFValues) {
    $cell = $ws.Cells.Item($item.Row, 32)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Text
}

# Column AG (33)
$AGValues = @(
    @{ Row = 2; Text = "`n" },
    @{ Row = 3; Text = " varifold`n" },
    @{ Row = 4; Text = " 10`n" },
    @{ Row = 5; Text = " 10`n" },
    @{ Row = 6; Text = " 0.05`n" },
    @{ Row = 7; Text = " 11`n" },
    @{ Row = 8; Text = " mask_05.vtk`n" },
    @{ Row = 9; Text = "`n" },
    @{ Row = 10; Text = " `n" },
    @{ Row = 11; Text = " 200.0`n" },
    @{ Row = 12; Text = " 1.5294790284563924`n" },
    @{ Row = 13; Text = " 1.0835121051448278`n" },
    @{ Row = 14; Text = " 0.2601479275466791`n" },
    @{ Row = 15; Text = " 0.8182787769669063`n" },
    @{ Row = 16; Text = " 1.1391222991510368`n" },
    @{ Row = 17; Text = " 1.9703321066752162`n" },
    @{ Row = 18; Text = " 7.040165472553342`n" },
    @{ Row = 19; Text = "`n" },
    @{ Row = 20; Text = " `n" },
    @{ Row = 21; Text = " 62.81028922143926`n" },
    @{ Row = 22; Text = " 53.00087159658319`n" },
    @{ Row = 23; Text = " 85.34230186030939`n" },
    @{ Row = 24; Text = " 59.64942912310182`n" },
    @{ Row = 25; Text = " 18.23559810113121`n" },
    @{ Row = 26; Text = " 55.80769798051297`n" },
    @{ Row = 27; Text = "`n" }
)
foreach ($item in $AGValues) {
    $cell = $ws.Cells.Item($item.Row, 33)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Text
}

